$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 517 ("「コーヒーは急には飲まれない」...") which shifts all subsequent rows up by one.
$ws.Rows.Item(517).Delete()
